# Std140_HE_Output.xlsx maintenance edit
# (commit: "Add reference to plotly and more HE support")
#
# Concrete content changes applied here:
#   1. Rename the single worksheet tab from "ESP-HOT" to "Sheet1".
#   2. Break/remove the stale external workbook link that pointed at
#      Std140_HE_Results.xlsx ("Informative Materials" results workbook).
#      No formula in the sheet actually references it, so breaking the
#      link drops both the <externalReferences> entry in workbook.xml and
#      the xl/externalLinks/externalLink1.xml part (+ its rels) entirely.

$wb = $excel.ActiveWorkbook

# --- 1) Rename worksheet tab -------------------------------------------
$ws = $wb.ActiveSheet
if ($ws.Name -eq "ESP-HOT") {
    $ws.Name = "Sheet1"
} else {
    # Fall back to renaming the first sheet if layout ever changes.
    $wb.Worksheets.Item(1).Name = "Sheet1"
}

# --- 2) Remove unused external link -------------------------------------
$links = $wb.LinkSources(1)
if ($links) {
    foreach ($link in $links) {
        $wb.BreakLink($link, 1)
    }
}

Write-Host ("Sheets: " + ($wb.Worksheets.Item(1).Name))
Write-Host ("Remaining link sources: " + ($wb.LinkSources(1)))
